$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "62.591.20"
Set-TextValue "E2" "  -0.70%  "

# Row 3
Set-TextValue "D3" "2.575.33"
Set-TextValue "E3" "  +1.04%  "

# Row 4
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  +0.06%  "

# Row 5
Set-TextValue "D5" "580.67"
Set-TextValue "E5" "  -0.22%  "

# Row 6
Set-TextValue "D6" "144.48"
Set-TextValue "E6" "  -1.71%  "

# Row 7
Set-TextValue "D7" "1.00"
Set-TextValue "E7" "  +0.06%  "

# Row 8
Set-TextValue "E8" "  +1.31%  "

# Row 9
Set-TextValue "D9" "0.106"
Set-TextValue "E9" "  +0.29%  "

# Row 10
Set-TextValue "E10" "  +0.22%  "

# Row 11
Set-TextValue "E11" "  -0.41%  "

# Row 12
Set-TextValue "E12" "  -0.67%  "

# Row 13
Set-TextValue "D13" "26.90"
Set-TextValue "E13" "  -2.31%  "

# Row 14
Set-TextValue "D14" "3.037.77"
Set-TextValue "E14" "  +1.16%  "

# Row 15
Set-TextValue "D15" "62.537.38"
Set-TextValue "E15" "  -0.62%  "

# Row 16
Set-TextValue "E16" "  +0.21%  "

# Row 17
Set-TextValue "D17" "2.580.27"
Set-TextValue "E17" "  +1.49%  "

# Row 18
Set-TextValue "E18" "  -1.20%  "

# Row 19
Set-TextValue "D19" "337.59"
Set-TextValue "E19" "  -0.19%  "

# Row 20
Set-TextValue "E20" "  +0.66%  "

# Row 21
Set-TextValue "D21" "6.63"
Set-TextValue "E21" "  -1.75%  "

# Row 22
Set-TextValue "D22" "0.998"
Set-TextValue "E22" "  -0.10%  "

# Row 23
Set-TextValue "D23" "67.01"
Set-TextValue "E23" "  +2.14%  "

# Row 24
Set-TextValue "D24" "2.702.31"
Set-TextValue "E24" "  +0.87%  "

# Row 25
Set-TextValue "D25" "0.166"
Set-TextValue "E25" "  -2.15%  "

# Row 26
Set-TextValue "E26" "  -1.90%  "

# Row 27
Set-TextValue "E27" "  -0.28%  "

# Row 28
Set-TextValue "B28" "SuiNetwork"
Set-TextValue "C28" "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue "D28" "1.47"
Set-TextValue "E28" "  -1.30%  "

# Row 29
Set-TextValue "B29" "Aptos"
Set-TextValue "C29" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D29" "7.86"
Set-TextValue "E29" "  +2.17%  "

# Row 30
Set-TextValue "D30" "8.19"
Set-TextValue "E30" "  -1.85%  "

# Row 31
Set-TextValue "E31" "  -1.72%  "

# Row 32
Set-TextValue "B32" "PEPE"
Set-TextValue "C32" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D32" "0.0₃0806"
Set-TextValue "E32" "  -1.11%  "

# Row 33
Set-TextValue "B33" "Bittensor"
Set-TextValue "C33" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D33" "459.14"
Set-TextValue "E33" "  +9.17%  "

# Row 34
Set-TextValue "D34" "176.55"
Set-TextValue "E34" "  -0.72%  "

# Row 35
Set-TextValue "E35" "  +2.98%  "

# Row 36
Set-TextValue "E36" "  +0.13%  "

# Row 37
Set-TextValue "D37" "0.399"
Set-TextValue "E37" "  -0.22%  "

# Row 38
Set-TextValue "D38" "18.85"
Set-TextValue "E38" "  -1.22%  "

# Row 39
Set-TextValue "D39" "4.45"
Set-TextValue "E39" "  +2.03%  "

# Row 40
Set-TextValue "E40" "  -0.02%  "

# Row 41
Set-TextValue "D41" "1.68"
Set-TextValue "E41" "  -3.35%  "

# Row 42
Set-TextValue "D42" "156.66"
Set-TextValue "E42" "  +4.05%  "

# Row 43
Set-TextValue "D43" "3.71"
Set-TextValue "E43" "  -1.86%  "

# Row 44
Set-TextValue "D44" "21.04"
Set-TextValue "E44" "  +1.50%  "

# Row 45
Set-TextValue "D45" "0.626"
Set-TextValue "E45" "  +3.86%  "

# Row 46
Set-TextValue "D46" "0.0534"
Set-TextValue "E46" "  -0.58%  "

# Row 47
Set-TextValue "D47" "0.0964"
Set-TextValue "E47" "  -0.39%  "

# Row 48
Set-TextValue "D48" "0.0233"
Set-TextValue "E48" "  -2.10%  "

# Row 49
Set-TextValue "D49" "18.02"
Set-TextValue "E49" "  -1.42%  "

# Row 50
Set-TextValue "E50" "  +0.94%  "

# Row 51
Set-TextValue "D51" "1.68"
Set-TextValue "E51" "  -1.33%  "
